$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary updates ---
$ws.Range("E11").Value = 1138800
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 4

# --- Expand table: insert 7 rows before row 29 so the table grows from 14 to 21 data rows ---
for ($i = 0; $i -lt 7; $i++) {
  $ws.Rows.Item(29).Insert()
}

# --- Copy the normal data-row formatting (row 28) onto the freshly inserted rows 29:35 ---
$ws.Range("B28:J28").Copy()
$ws.Range("B29:J35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Write the refreshed worker/period data set (rows 16-36) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73161753"
$ws.Range("D16").Value = "ORLANDO MANUEL MERCADO DIAZ"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "18881781"
$ws.Range("D17").Value = "JULIO RAFAEL RIVERO OLIVERA"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1102148287"
$ws.Range("D18").Value = "TATIANA LUCIA MONTES BENITEZ"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73161753"
$ws.Range("D19").Value = "ORLANDO MANUEL MERCADO DIAZ"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "18881781"
$ws.Range("D20").Value = "JULIO RAFAEL RIVERO OLIVERA"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1102148287"
$ws.Range("D21").Value = "TATIANA LUCIA MONTES BENITEZ"
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1052950363"
$ws.Range("D22").Value = "MAYLEN MARCELA MENCO HERNANDEZ"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73161753"
$ws.Range("D23").Value = "ORLANDO MANUEL MERCADO DIAZ"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "18881781"
$ws.Range("D24").Value = "JULIO RAFAEL RIVERO OLIVERA"
$ws.Range("E24").Value = "2507"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1102148287"
$ws.Range("D25").Value = "TATIANA LUCIA MONTES BENITEZ"
$ws.Range("E25").Value = "2507"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1085180410"
$ws.Range("D26").Value = "LUIS FERNANDO RAMOS GOMEZ"
$ws.Range("E26").Value = "2507"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 1423500
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1010113650"
$ws.Range("D27").Value = "WENDY TATIANA TABARES PORTELA"
$ws.Range("E27").Value = "2507"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1052950363"
$ws.Range("D28").Value = "MAYLEN MARCELA MENCO HERNANDEZ"
$ws.Range("E28").Value = "2508"
$ws.Range("F28").Value = 56940
$ws.Range("G28").Value = 1423500
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "73161753"
$ws.Range("D29").Value = "ORLANDO MANUEL MERCADO DIAZ"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 56940
$ws.Range("G29").Value = 1423500
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "18881781"
$ws.Range("D30").Value = "JULIO RAFAEL RIVERO OLIVERA"
$ws.Range("E30").Value = "2508"
$ws.Range("F30").Value = 56940
$ws.Range("G30").Value = 1423500
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1102148287"
$ws.Range("D31").Value = "TATIANA LUCIA MONTES BENITEZ"
$ws.Range("E31").Value = "2508"
$ws.Range("F31").Value = 56940
$ws.Range("G31").Value = 1423500
$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1085180410"
$ws.Range("D32").Value = "LUIS FERNANDO RAMOS GOMEZ"
$ws.Range("E32").Value = "2508"
$ws.Range("F32").Value = 56940
$ws.Range("G32").Value = 1423500
$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1010113650"
$ws.Range("D33").Value = "WENDY TATIANA TABARES PORTELA"
$ws.Range("E33").Value = "2508"
$ws.Range("F33").Value = 56940
$ws.Range("G33").Value = 1423500
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "9023928"
$ws.Range("D34").Value = "CARLOS ENRIQUE BENITEZ PIÃ?ERES"
$ws.Range("E34").Value = "2508"
$ws.Range("F34").Value = 37960
$ws.Range("G34").Value = 1423500
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1006003800"
$ws.Range("D35").Value = "JOSE MANUEL RODRIGUEZ GUZMAN"
$ws.Range("E35").Value = "2508"
$ws.Range("F35").Value = 37960
$ws.Range("G35").Value = 1423500
$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1001834172"
$ws.Range("D36").Value = "JOSE GREGORIO BLANCO MEJIA"
$ws.Range("E36").Value = "2508"
$ws.Range("F36").Value = 37960
$ws.Range("G36").Value = 1423500
